$d = $word.ActiveDocument

# Remove the first three paragraphs (Title, Author, Date)
$d.Paragraphs(1).Range.Delete()
$d.Paragraphs(1).Range.Delete()
$d.Paragraphs(1).Range.Delete()

# Fix capitalization: "markdown" -> "Markdown"
$d.Content.Find.Execute("This is a markdown file", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This is a Markdown file", 2)
